# Natmi following Dr Hou advice
# Rebuild the LR-pair table (Cxcl13 -> Ackr4) with the full 3x3 cross of
# Sending cluster x Target cluster ("FAPs", "M2", "sCs"), replacing the
# previous 2-row (FAPs/FAPs x ECs/FAPs) table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one data row: A..T in order
#   A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
#   E..T the 16 numeric metric columns.
$rowsData = @(
    @("FAPs","Cxcl13","Ackr4","ECs",3,1,4.046465333333333,12.139396,0.8279246837497715,0.8279246837497715,3,1,0.2315733333333333,0.69472,0.2371078251520917,0.2371078251520917,0.9370534654577777,8.43348118912,0.1963074211536416,0.1963074211536416),
    @("FAPs","Cxcl13","Ackr4","FAPs",3,1,4.046465333333333,12.139396,0.8279246837497715,0.8279246837497715,3,1,0.7200953333333334,2.160286,0.7373052671097876,0.7373052671097876,2.913840803028445,26.224567227256,0.6104332300989117,0.6104332300989117),
    @("FAPs","Cxcl13","Ackr4","sCs",3,1,4.046465333333333,12.139396,0.8279246837497715,0.8279246837497715,1,0.3333333333333333,0.02498966666666666,0.07496899999999999,0.02558690773812063,0.02558690773812063,0.1011198198582222,0.9100783787239999,0.0211840324972181,0.0211840324972181),
    @("M2","Cxcl13","Ackr4","ECs",2,0.6666666666666666,0.2754046666666667,0.826214,0.05634901148785605,0.05634901148785604,3,1,0.2315733333333333,0.69472,0.2371078251520917,0.2371078251520917,0.06377637667555555,0.57398739008,0.01336079156335578,0.01336079156335578),
    @("M2","Cxcl13","Ackr4","FAPs",2,0.6666666666666666,0.2754046666666667,0.826214,0.05634901148785605,0.05634901148785604,3,1,0.7200953333333334,2.160286,0.7373052671097876,0.7373052671097876,0.1983176152448889,1.784858537204,0.0415464229664262,0.04154642296642619),
    @("M2","Cxcl13","Ackr4","sCs",2,0.6666666666666666,0.2754046666666667,0.826214,0.05634901148785605,0.05634901148785604,1,0.3333333333333333,0.02498966666666666,0.07496899999999999,0.02558690773812063,0.02558690773812063,0.006882270818444444,0.061940437366,0.001441796958074072,0.001441796958074072),
    @("sCs","Cxcl13","Ackr4","ECs",2,0.6666666666666666,0.5656100000000001,1.69683,0.1157263047623724,0.1157263047623724,3,1,0.2315733333333333,0.69472,0.2371078251520917,0.2371078251520917,0.1309801930666667,1.1788217376,0.02743961243509429,0.02743961243509428),
    @("sCs","Cxcl13","Ackr4","FAPs",2,0.6666666666666666,0.5656100000000001,1.69683,0.1157263047623724,0.1157263047623724,3,1,0.7200953333333334,2.160286,0.7373052671097876,0.7373052671097876,0.4072931214866667,3.665638093380001,0.0853256140444497,0.08532561404444969),
    @("sCs","Cxcl13","Ackr4","sCs",2,0.6666666666666666,0.5656100000000001,1.69683,0.1157263047623724,0.1157263047623724,1,0.3333333333333333,0.02498966666666666,0.07496899999999999,0.02558690773812063,0.02558690773812063,0.01413440536333333,0.12720964827,0.002961078282828453,0.002961078282828454)
)

$startRow = 2
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $rowsData[$i]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
